# Update column L values: "rtd:appendix-1/" -> "rtd:appendices/appendix-1/"
# Applies to the data rows that reference the rtd appendix vocabulary (rows 18-146).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 18; $row -le 146; $row++) {
    $cell = $ws.Cells.Item($row, 12)   # column L = 12
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("rtd:appendix-1/")) {
        $cell.Value2 = $val.ToString().Replace("rtd:appendix-1/", "rtd:appendices/appendix-1/")
    }
}
